# Scheduled runner update: refresh market price / profit figures on each
# crafting-job Leve sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2745.762
$ws.Range("I135").Value = 1448.8125
$ws.Range("J135").Value = 6896
$ws.Range("K135").Value = 13039.3125
$ws.Range("L135").Value = 62064
$ws.Range("M135").Value = -10504.3125
$ws.Range("N135").Value = -67134

$ws.Range("H137").Value = 34052.027
$ws.Range("J137").Value = 5418.8
$ws.Range("L137").Value = 16256.4
$ws.Range("N137").Value = -21356.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6777.857
$ws.Range("I2").Value = 8964.143
$ws.Range("K2").Value = 8964.143
$ws.Range("M2").Value = -8851.143

$ws.Range("H32").Value = 8279.485000000001
$ws.Range("I32").Value = 7844.4204
$ws.Range("K32").Value = 7844.4204
$ws.Range("M32").Value = -7557.4204

$ws.Range("H61").Value = 10754.538
$ws.Range("I61").Value = 3904.4
$ws.Range("J61").Value = 15035.875
$ws.Range("K61").Value = 3904.4
$ws.Range("L61").Value = 15035.875
$ws.Range("M61").Value = -3692.4
$ws.Range("N61").Value = -15459.875

$ws.Range("H110").Value = 10863.2
$ws.Range("I110").Value = 10610
$ws.Range("K110").Value = 10610
$ws.Range("M110").Value = -8565

$ws.Range("H116").Value = 6777.857
$ws.Range("I116").Value = 8964.143
$ws.Range("K116").Value = 8964.143
$ws.Range("M116").Value = -6670.143

$ws.Range("H136").Value = 10754.538
$ws.Range("I136").Value = 3904.4
$ws.Range("J136").Value = 15035.875
$ws.Range("K136").Value = 11713.2
$ws.Range("L136").Value = 45107.625
$ws.Range("M136").Value = -9163.200000000001
$ws.Range("N136").Value = -50207.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6777.857
$ws.Range("I3").Value = 8964.143
$ws.Range("K3").Value = 8964.143
$ws.Range("M3").Value = -8850.143

$ws.Range("H134").Value = 1839.8182
$ws.Range("I134").Value = 1331.4
$ws.Range("J134").Value = 3428.625
$ws.Range("K134").Value = 3994.2
$ws.Range("L134").Value = 10285.875
$ws.Range("M134").Value = -1459.2
$ws.Range("N134").Value = -15355.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 158114.12
$ws.Range("I31").Value = 206007.78
$ws.Range("K31").Value = 206007.78
$ws.Range("M31").Value = -205712.78

$ws.Range("H34").Value = 158114.12
$ws.Range("I34").Value = 206007.78
$ws.Range("K34").Value = 206007.78
$ws.Range("M34").Value = -205805.78

$ws.Range("H58").Value = 2980.818
$ws.Range("I58").Value = 2532.111
$ws.Range("K58").Value = 2532.111
$ws.Range("M58").Value = -2329.111

$ws.Range("H70").Value = 60000
$ws.Range("J70").Value = 60000
$ws.Range("L70").Value = 60000
$ws.Range("N70").Value = -60630

$ws.Range("H73").Value = 60000
$ws.Range("J73").Value = 60000
$ws.Range("L73").Value = 60000
$ws.Range("N73").Value = -62184

$ws.Range("H107").Value = 3584.6296
$ws.Range("I107").Value = 617.64703
$ws.Range("K107").Value = 617.64703
$ws.Range("M107").Value = 1302.35297

$ws.Range("H132").Value = 5501.8237
$ws.Range("I132").Value = 2720.6875
$ws.Range("K132").Value = 8162.0625
$ws.Range("M132").Value = -5632.0625

$ws.Range("H136").Value = 2980.818
$ws.Range("I136").Value = 2532.111
$ws.Range("K136").Value = 7596.333
$ws.Range("M136").Value = -5046.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 726.25
$ws.Range("I18").Value = 535.5
$ws.Range("K18").Value = 1606.5
$ws.Range("M18").Value = -1437.5

$ws.Range("H47").Value = 13435.941
$ws.Range("I47").Value = 9998.6
$ws.Range("J47").Value = 18346.428
$ws.Range("K47").Value = 29995.8
$ws.Range("L47").Value = 55039.284
$ws.Range("M47").Value = -29564.8
$ws.Range("N47").Value = -55901.284

$ws.Range("H49").Value = 850
$ws.Range("I49").Value = 850
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2550
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("M49").Value = -2394

$ws.Range("H50").Value = 427
$ws.Range("I50").Value = 443.6
$ws.Range("K50").Value = 1330.8
$ws.Range("M50").Value = -849.8000000000002

$ws.Range("H53").Value = 427
$ws.Range("I53").Value = 443.6
$ws.Range("K53").Value = 1330.8
$ws.Range("M53").Value = -849.8000000000002

$ws.Range("H55").Value = 6379.25
$ws.Range("J55").Value = 6718.1665
$ws.Range("L55").Value = 20154.4995
$ws.Range("N55").Value = -20508.4995

$ws.Range("H131").Value = 7043787
$ws.Range("J131").Value = 1598.9219
$ws.Range("L131").Value = 4796.7657
$ws.Range("N131").Value = -14876.7657

$ws.Range("H140").Value = 5320931
$ws.Range("I140").Value = 25001100
$ws.Range("K140").Value = 75003300
$ws.Range("M140").Value = -74998120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1166
$ws.Range("I31").Value = 1166
$ws.Range("K31").Value = 1166
$ws.Range("M31").Value = -874

$ws.Range("H37").Value = 1166
$ws.Range("I37").Value = 1166
$ws.Range("K37").Value = 1166
$ws.Range("M37").Value = -889

$ws.Range("H132").Value = 107136.63
$ws.Range("I132").Value = 140434.5
$ws.Range("J132").Value = 18342.334
$ws.Range("K132").Value = 421303.5
$ws.Range("L132").Value = 55027.00199999999
$ws.Range("M132").Value = -418773.5
$ws.Range("N132").Value = -60087.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 549.7826
$ws.Range("I16").Value = 549.7826
$ws.Range("K16").Value = 549.7826
$ws.Range("M16").Value = -379.7826

$ws.Range("H55").Value = 6794.8335
$ws.Range("I55").Value = 543.6
$ws.Range("K55").Value = 543.6
$ws.Range("M55").Value = -370.6

$ws.Range("H64").Value = 34283.5
$ws.Range("J64").Value = 34283.5
$ws.Range("L64").Value = 34283.5
$ws.Range("N64").Value = -34733.5

$ws.Range("H67").Value = 34283.5
$ws.Range("J67").Value = 34283.5
$ws.Range("L67").Value = 34283.5
$ws.Range("N67").Value = -35843.5

$ws.Range("H132").Value = 12619.92
$ws.Range("I132").Value = 17556.4
$ws.Range("K132").Value = 52669.2
$ws.Range("M132").Value = -50139.2

$ws.Range("H136").Value = 7506.25
$ws.Range("I136").Value = 3643.182
$ws.Range("K136").Value = 10929.546
$ws.Range("M136").Value = -8379.545999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 25000
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 25000
$ws.Range("N63").Value = -26248

$ws.Range("H66").Value = 25000
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 75000
$ws.Range("N66").Value = -81240

$ws.Range("H132").Value = 979.8
$ws.Range("I132").Value = 909.8333
$ws.Range("J132").Value = 1259.6666
$ws.Range("K132").Value = 2729.4999
$ws.Range("L132").Value = 3778.9998
$ws.Range("M132").Value = -199.4998999999998
$ws.Range("N132").Value = -8838.9998
